$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 256-269 (daily figures through 27/05) to append below existing data (which ends at row 255).
$data = @(
  @(44330, 3, 13, 131.8191036300953),
  @(44331, 2, 12, 121.6791725816264),
  @(44332, 3, 14, 141.9590346785642),
  @(44333, 2, 13, 131.8191036300953),
  @(44334, 7, 19, 192.6586899209085),
  @(44335, 1, 19, 192.6586899209085),
  @(44336, 2, 20, 202.7986209693774),
  @(44337, 1, 18, 182.5187588724397),
  @(44338, 0, 16, 162.2388967755019),
  @(44339, 0, 13, 131.8191036300953),
  @(44340, 2, 13, 131.8191036300953),
  @(44341, 1, 7, 70.9795173392821),
  @(44342, 0, 6, 60.83958629081322),
  @(44343, 0, 4, 40.55972419387548)
)

$startRow = 256
$endRow = $startRow + $data.Count - 1

# Propagate the date-column formatting (bold, bordered, centered, custom date format)
# from the last existing row down across all the new rows before filling in values.
$ws.Range("A255").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  $ws.Range("A$r").Value = $row[0]
  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
}
